# Insert a new data row at row 98 (pushing the existing rows 98-177 down to
# 99-178) and populate the new row with the new price-observation record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("98:98").Insert()

$ws.Range("A98").Value = 10
$ws.Range("B98").Value = "Vega Modelo de Temuco"
$ws.Range("C98").Value = "La Araucanía"
$ws.Range("D98").Value = 44827
$ws.Range("E98").Value = 9
$ws.Range("F98").Value = "Fruta"
$ws.Range("G98").Value = 100104
$ws.Range("H98").Value = "Frutos de pepita"
$ws.Range("I98").Value = 100104001
$ws.Range("J98").Value = "Granada"
$ws.Range("K98").Value = "Wonderfull"
$ws.Range("L98").Value = "Primera"
$ws.Range("M98").Value = 55
$ws.Range("N98").Value = 14000
$ws.Range("O98").Value = 14000
$ws.Range("P98").Value = 14000
$ws.Range("Q98").Value = "`$/bandeja 10 kilos granel"
$ws.Range("R98").Value = "Provincia de Limarí"
$ws.Range("S98").Value = 1400
$ws.Range("T98").Value = 10
